$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 33: mark "83. Remove Duplicates from Sorted List" with an orange
#     highlight and fix up the Pattern column (was "?", now "Fast and slow
#     pointers"), plus a newly-touched (empty) F33 cell. ---
$orange = 49407   # RGB(255,192,0) == FFC000, packed BGR for OLE Color

$ws.Range("D33").Value = "Fast and slow pointers"

$ws.Range("C33").Style = "Гиперссылка"

$ws.Range("A33:F33").Interior.Color = $orange

# --- New rows appended at the bottom of the table ---
$ws.Range("A66").Formula = "=ROW()-1"
$ws.Range("B66").Value = "1971. Find if Path Exists in Graph"
$ws.Range("C66").Value = "https://leetcode.com/problems/find-if-path-exists-in-graph/"
$ws.Range("C66").Style = "Гиперссылка"
$ws.Range("D66").Value = "BFS"
$ws.Range("E66").Value = "Easy"

$ws.Range("A67").Value = 65
$ws.Range("B67").Value = "261. Graph Valid Tree"
$ws.Range("C67").Value = "https://leetcode.com/problems/graph-valid-tree/"
$ws.Range("C67").Style = "Гиперссылка"
$ws.Range("D67").Value = "Union-find"
$ws.Range("E67").Value = "Medium"

# --- Update view / selection to where the new rows were edited ---
$null = $ws.Range("B46").Select()
